$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, copying the
# existing header style (bold/border/centered, same as O1) and setting
# the next sequential values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For the data rows (2-25): swap the I/K values and the M/O values,
# then append the two new columns P and Q with value 2.
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
